# Cleanup pass: remove the decorative separator lines, the screenshot
# images, and the leftover empty "spacing before" paragraphs that used to
# pad the document around code-sample tables.
#
# Paragraph indices (1-based, as in the *original* document) of every
# paragraph that must disappear completely (including its paragraph
# mark), from first to last:
#   2   - screenshot image under the title
#   5   - screenshot image under "Introducción"
#   8   - "────" separator after the intro text
#   36  - empty spacer paragraph after the HTML code table
#   63  - empty spacer paragraph after the vocabulary code table
#   93  - empty spacer paragraph after the generation-logic code table
#   116 - empty spacer paragraph after the capitalisation code table
#   132 - empty spacer paragraph after the CSS code table
#   157 - empty spacer paragraph after the border-color code table
#   177 - empty spacer paragraph after the accessibility code table
#   178 - "────" separator before "Presentación del proyecto"
#   180 - screenshot image under "Presentación del proyecto"
#   185 - "────" separator before "Conclusión"
#
# Deleting a paragraph's Range removes its paragraph mark too, so later
# paragraphs shift down by one. Walking the index list from the highest
# number to the lowest keeps every not-yet-processed index valid.

$d = $word.ActiveDocument

$targets = @(185, 180, 178, 177, 157, 132, 116, 93, 63, 36, 8, 5, 2)

foreach ($idx in $targets) {
    $p = $d.Paragraphs.Item($idx)
    $p.Range.Delete()
}
